# Insert a new weekly price record as row 136 on the single data sheet,
# pushing the existing rows 136:210 down to 137:211 (dimension grows to
# A1:R211). The new row carries a fresh "Perejil" (parsley) price quote.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 136:210 down one row to make room for the new record.
$ws.Rows.Item(136).Insert()

# Populate the newly inserted row 136 with the new data point.
$ws.Range("A136").Value = 9
$ws.Range("B136").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C136").Value = 'Metropolitana'
$ws.Range("D136").Value = 44460
$ws.Range("E136").Value = 13
$ws.Range("F136").Value = 100112044
$ws.Range("G136").Value = 'Perejil'
$ws.Range("H136").Value = 'Sin especificar'
$ws.Range("I136").Value = 'Primera'
$ws.Range("J136").Value = 97
$ws.Range("K136").Value = 9000
$ws.Range("L136").Value = 10000
$ws.Range("M136").Value = 9495
$ws.Range("N136").Value = '$/docena de atados'
$ws.Range("O136").Value = 'Región Metropolitana'
$ws.Range("P136").Value = 3165
$ws.Range("Q136").Value = 3
$ws.Range("R136").Value = 'Hortaliza'
